$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 5

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 3

# Add new rows 4 and 5, matching the style of A2:A3 (bold/bordered/centered) for column A
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 2

# Copy style from A3 (style index 1) to A4:A5 so they match the existing formatted column A cells
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
